$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the diff: (cell -> new value).
# Numeric-looking text values are written via a text-format/clear-format
# round-trip so Excel stores them as literal text (matching the source
# inlineStr cells) instead of silently coercing them to numbers, while
# leaving the cell's style index untouched (index 0, same as before).

$textValues = @{
    D2 = '69.504.97'
    E2 = '  +0.27%  '
    D3 = '2.493.62'
    E3 = '  -0.96%  '
    E4 = '  +0.06%  '
    E5 = '  -0.41%  '
    E6 = '  +0.54%  '
    E7 = '  +0.02%  '
    E8 = '  -1.17%  '
    E9 = '  -0.45%  '
    E10 = '  -0.74%  '
    E11 = '  -1.58%  '
    E12 = '  -0.70%  '
    D13 = '2.951.72'
    E13 = '  -1.05%  '
    D14 = '69.434.94'
    E14 = '  +0.37%  '
    E15 = '  -0.17%  '
    E16 = '  -2.32%  '
    D17 = '2.462.38'
    E17 = '  -2.32%  '
    E18 = '  -0.73%  '
    E19 = '  +1.90%  '
    E20 = '  -3.51%  '
    E21 = '  -0.22%  '
    E22 = '  -4.16%  '
    E23 = '  -0.03%  '
    E24 = '  -1.58%  '
    E25 = '  -3.75%  '
    D26 = '2.623.62'
    E26 = '  -0.96%  '
    E27 = '  -2.59%  '
    E28 = '  -0.08%  '
    D29 = '0.0₃0874'
    E29 = '  -1.78%  '
    E30 = '  -2.57%  '
    E31 = '  -4.21%  '
    E32 = '  -3.24%  '
    E33 = '  +0.06%  '
    E34 = '  -1.04%  '
    E35 = '  +96.77%  '
    B36 = 'Kaspa'
    C36 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    E36 = '  -2.91%  '
    B37 = 'Monero'
    C37 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    E37 = '  -2.79%  '
    E39 = '  -1.96%  '
    E40 = '  -0.02%  '
    E41 = '  -1.34%  '
    E42 = '  -2.23%  '
    E43 = '  -1.35%  '
    E44 = '  -1.46%  '
    E45 = '  -4.90%  '
    E46 = '  -2.24%  '
    E47 = '  -1.08%  '
    E48 = '  -3.14%  '
    E49 = '  -1.14%  '
    E50 = '  -1.45%  '
    E51 = '  -0.56%  '
}

$numericTextValues = @{
    D5 = '570.03'
    D6 = '166.47'
    D8 = '0.509'
    D9 = '0.159'
    D16 = '24.23'
    D18 = '11.24'
    D19 = '354.24'
    D20 = '7.38'
    D24 = '69.38'
    D25 = '3.80'
    D27 = '8.63'
    D30 = '7.61'
    D31 = '441.14'
    D35 = '3.07'
    D36 = '0.113'
    D37 = '152.96'
    D39 = '18.16'
    D42 = '4.59'
    D44 = '2.19'
    D46 = '138.63'
}

foreach ($cell in $textValues.Keys) {
    $ws.Range($cell).Value = $textValues[$cell]
}

foreach ($cell in $numericTextValues.Keys) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $numericTextValues[$cell]
    $r.ClearFormats()
}

